$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended at the bottom of the sheet (rows 245-247),
# continuing the existing daily series in columns A-D.
$newRows = @(
    @{ Row = 245; Date = 44319; B = 0; C = 13; D = 85.77461071522829 },
    @{ Row = 246; Date = 44320; B = 1; C = 13; D = 85.77461071522829 },
    @{ Row = 247; Date = 44321; B = 0; C = 12; D = 79.1765637371338 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

# Column A carries a date-formatted style (same as the rest of the series);
# copy that formatting from the last previously-existing row (244) onto the
# newly added date cells so the new cells share style index "s=2".
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)
$excel.CutCopyMode = $false
